$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K
$ws.Range("K2").Value = "test_xlr_n_percent"

# Build the target style (font1, numFmt General, right/bottom align) on a scratch cell
$ws.Range("G3").Copy()
$scratch = $ws.Range("Z100")
$scratch.PasteSpecial(-4122)
$scratch.HorizontalAlignment = -4152
$scratch.VerticalAlignment = -4107

# Apply that resolved style to the whole K3:K35 range in one shot
$scratch.Copy()
$ws.Range("K3:K35").PasteSpecial(-4122)
$scratch.Clear()

# Fill in the values for K3:K34 (n (pct%) strings); K35 stays blank
$ws.Cells.Item(3, 11).Value = "1 (3%)"
$ws.Cells.Item(4, 11).Value = "2 (6%)"
$ws.Cells.Item(5, 11).Value = "3 (9%)"
$ws.Cells.Item(6, 11).Value = "4 (12%)"
$ws.Cells.Item(7, 11).Value = "5 (16%)"
$ws.Cells.Item(8, 11).Value = "6 (19%)"
$ws.Cells.Item(9, 11).Value = "7 (22%)"
$ws.Cells.Item(10, 11).Value = "8 (25%)"
$ws.Cells.Item(11, 11).Value = "9 (28%)"
$ws.Cells.Item(12, 11).Value = "10 (31%)"
$ws.Cells.Item(13, 11).Value = "11 (34%)"
$ws.Cells.Item(14, 11).Value = "12 (38%)"
$ws.Cells.Item(15, 11).Value = "13 (41%)"
$ws.Cells.Item(16, 11).Value = "14 (44%)"
$ws.Cells.Item(17, 11).Value = "15 (47%)"
$ws.Cells.Item(18, 11).Value = "16 (50%)"
$ws.Cells.Item(19, 11).Value = "17 (53%)"
$ws.Cells.Item(20, 11).Value = "18 (56%)"
$ws.Cells.Item(21, 11).Value = "19 (59%)"
$ws.Cells.Item(22, 11).Value = "20 (62%)"
$ws.Cells.Item(23, 11).Value = "21 (66%)"
$ws.Cells.Item(24, 11).Value = "22 (69%)"
$ws.Cells.Item(25, 11).Value = "23 (72%)"
$ws.Cells.Item(26, 11).Value = "24 (75%)"
$ws.Cells.Item(27, 11).Value = "25 (78%)"
$ws.Cells.Item(28, 11).Value = "26 (81%)"
$ws.Cells.Item(29, 11).Value = "27 (84%)"
$ws.Cells.Item(30, 11).Value = "28 (88%)"
$ws.Cells.Item(31, 11).Value = "29 (91%)"
$ws.Cells.Item(32, 11).Value = "30 (94%)"
$ws.Cells.Item(33, 11).Value = "31 (97%)"
$ws.Cells.Item(34, 11).Value = "32 (100%)"
